$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update coin price / volume(1h) cells to match the refreshed symbol-list snapshot.
# Each target cell is stored as text (e.g. "243.80", "-0.73%"), so the number format
# is pinned to Text ("@") before writing the value to prevent Excel from auto-coercing
# the numeric-looking string into a real number/percentage.
$updates = @{
    "D2" = "243.80"
    "E2" = "-0.73%"
    "D3" = "26.63"
    "E3" = "4.56%"
    "D4" = "5.152"
    "E4" = "0.42%"
    "D5" = "0.05617"
    "E5" = "0.45%"
    "D6" = "6.493"
    "E6" = "0.12%"
    "D7" = "0.8187"
    "E7" = "0.00%"
    "D8" = "0.8321"
    "E8" = "-1.76%"
    "D9" = "0.1328"
    "E9" = "-1.03%"
    "D10" = "0.06923"
    "E10" = "-0.51%"
    "D11" = "0.02897"
    "E11" = "1.26%"
    "D12" = "0.09384"
    "E12" = "-0.04%"
    "D13" = "0.001509"
    "E13" = "0.09%"
    "D14" = "0.0006009"
    "E14" = "1.17%"
    "D15" = "0.006249"
    "E15" = "2.13%"
    "D16" = "3.647"
    "E16" = "3.29%"
    "E17" = "-0.09%"
    "E18" = "8.64%"
    "D20" = "0.03104"
    "E20" = "-4.12%"
    "D21" = "0.1292"
    "D22" = "3.740"
    "E22" = "-0.05%"
    "D23" = "0.04590"
    "E23" = "-1.81%"
    "E24" = "-2.47%"
    "E25" = "-1.73%"
    "D26" = "0.004491"
    "E26" = "-2.43%"
    "D27" = "0.00009800"
    "E27" = "2.09%"
    "E28" = "0.68%"
    "D40" = "0.03636"
    "E40" = "-0.52%"
    "D41" = "0.006180"
    "E41" = "1.15%"
    "D42" = "0.1053"
    "E42" = "0.08%"
    "D43" = "0.002570"
    "E43" = "4.26%"
    "D44" = "0.008098"
    "E44" = "3.98%"
    "D45" = "0.00005351"
    "E45" = "0.62%"
    "E46" = "0.00%"
    "D48" = "0.002573"
    "E48" = "20.97%"
    "E49" = "0.00%"
    "E50" = "0.00%"
}

foreach ($cellRef in $updates.Keys) {
    $cell = $ws.Range($cellRef)
    $cell.NumberFormat = "@"
    $cell.Value = $updates[$cellRef]
}
